$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.685.77'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.965.66'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '244.18'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  +0.56%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '58.43'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('E8').Value = '  -0.01%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.372'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.36%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0804'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -3.48%  '
$ws.Range('E11').Value = '  -0.04%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '22.15'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +3.31%  '
$ws.Range('D13').Value = '2.256.91'
$ws.Range('E13').Value = '  +1.44%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.822'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.22%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '13.67'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +1.03%  '
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '1.960.62'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').Value = '36.616.74'
$ws.Range('E18').Value = '  +0.43%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '69.70'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('E21').Value = '  +1.95%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '228.51'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('E25').Value = '  +2.14%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '9.41'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +9.72%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '160.35'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -1.33%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '19.36'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  +1.23%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.13'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -2.06%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.70'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('E35').Value = '  -0.02%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '6.09'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.16%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '3.42'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +16.16%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +3.99%  '
$ws.Range('E39').Value = '  -0.66%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.0997'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +2.76%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.0212'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.16'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -0.73%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '16.01'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').Value = '1.368.86'
$ws.Range('E45').Value = '  +1.51%  '
$ws.Range('E46').Value = '  +0.32%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '87.40'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -0.21%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('D50').Value = '2.147.97'
$ws.Range('E50').Value = '  +1.66%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '43.45'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -4.86%  '
